# Test Report.xlsx edit: populate ContentProviderTests results, resize
# columns, move the selection, and set print/page setup to match the
# author's commit ("Updated Test Report and comment in
# ContentProviderTests.java").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test-result rows (Content Provider / SQLite test block) ---------

# "Content Provider (SQLite) tests:" section header
$ws.Range("A3").Value = "Content Provider (SQLite) tests:"

# testInserts() sub-tests 1-5
$ws.Range("A4").Value = "testInserts() sub-test 1"
$ws.Range("B4").Value = "PASSED! Returned no exceptions."
$ws.Range("C4").Value = "Inserts one Drink into the table: ""DrinkTable"" by using the Content Provider."

$ws.Range("A5").Value = "testInserts() sub-test 2"
$ws.Range("B5").Value = "PASSED! Returned one row from table."
$ws.Range("C5").Value = """DrinkTable"" should contain only one Drink."

$ws.Range("A6").Value = "testInserts() sub-test 3"
$ws.Range("B6").Value = "PASSED! Returned no exceptions."
$ws.Range("C6").Value = "Moves to the first record with cursor. Should point to the first and only row."

$ws.Range("A7").Value = "testInserts() sub-test 4"
$ws.Range("B7").Value = "PASSED! Returned correct data."
$ws.Range("C7").Value = "Check if all data in the entry is correct."

$ws.Range("A8").Value = "testInserts() sub-test 5"
$ws.Range("B8").Value = "PASSED! Returned an exception."
$ws.Range("C8").Value = "Try to insert a row that already exists. Should return an exception."

# testDeletes() sub-tests 1-3
$ws.Range("A10").Value = "testDeletes() sub-test 1"
$ws.Range("B10").Value = "PASSED! Returned zero rows deleted."
$ws.Range("C10").Value = "Try to delete an empty record."

$ws.Range("A11").Value = "testDeletes() sub-test 2"
$ws.Range("B11").Value = "PASSED! Returned one row deleted."
$ws.Range("C11").Value = "Delete an existing record."

$ws.Range("A12").Value = "testDeletes() sub-test 3"
$ws.Range("B12").Value = "PASSED! Query and check that cursor.getCount returns zero."
$ws.Range("C12").Value = "Check if the row was actually deleted."

# testUpdates() sub-tests 1-2
$ws.Range("A14").Value = "testUpdates() sub-test 1"
$ws.Range("B14").Value = "PASSED! Returned zero rows updated."
$ws.Range("C14").Value = "Try to update an empty record."

$ws.Range("A15").Value = "testUpdates() sub-test 2"
$ws.Range("B15").Value = "PASSED! Returned one row updated."
$ws.Range("C15").Value = "Update an existing record."

# --- Column widths (widened slightly) -------------------------------------
# ColumnWidth is character-width based and gets snapped to the host's pixel
# grid, so feed it the offset that lands closest to the target XML width.
$ws.Columns.Item(1).ColumnWidth = 28.7369791666667
$ws.Columns.Item(2).ColumnWidth = 55.0221354166667
$ws.Columns.Item(3).ColumnWidth = 67.7369791666667

# --- Selection moves from C2 to A2 -----------------------------------------
$ws.Range("A2").Select() | Out-Null

# --- Page setup: A4 portrait (adds <pageSetup .../> on save) --------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
